$d = $word.ActiveDocument

function New-WordOpenXml($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Change 1: expand the "Boot your Linux system..." instruction text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Boot your Linux system or VM, log in, and then open a terminal window and start the lab:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A. Boot your Linux system or VM.  If necessary, log in and then open a terminal window and cd to the labtainer/labtainer-student directory.  The pre-packaged Labtainer VM will start with such a terminal open for you.   Then start the lab:",
    2)

# ---------------------------------------------------------------------------
# Change 2: "cd labtainer/labtainer-student" paragraph -> "labtainer pubkey"
# (split into two runs) and clear the paragraph mark's run properties.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$tabsXml = '<w:tab w:val="left" w:pos="560" w:leader="none"/><w:tab w:val="left" w:pos="1120" w:leader="none"/><w:tab w:val="left" w:pos="1680" w:leader="none"/><w:tab w:val="left" w:pos="2240" w:leader="none"/><w:tab w:val="left" w:pos="2800" w:leader="none"/><w:tab w:val="left" w:pos="3360" w:leader="none"/><w:tab w:val="left" w:pos="3920" w:leader="none"/><w:tab w:val="left" w:pos="4480" w:leader="none"/><w:tab w:val="left" w:pos="5040" w:leader="none"/><w:tab w:val="left" w:pos="5600" w:leader="none"/><w:tab w:val="left" w:pos="6160" w:leader="none"/><w:tab w:val="left" w:pos="6720" w:leader="none"/>'
$runRPr = '<w:rFonts w:cs="Courier New" w:ascii="Courier New" w:hAnsi="Courier New"/><w:color w:val="000000"/>'
$p6Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:widowControl w:val="false"/><w:tabs>' + $tabsXml + '</w:tabs><w:rPr/></w:pPr><w:r><w:rPr>' + $runRPr + '</w:rPr><w:tab/><w:t>labtainer pubke</w:t></w:r><w:r><w:rPr>' + $runRPr + '</w:rPr><w:t>y</w:t></w:r></w:p></w:body>'
$p6.Range.InsertXML((New-WordOpenXml $p6Body))

# ---------------------------------------------------------------------------
# Change 3: "./start.py pubkey" paragraph -> empty run, keep paragraph
# properties (tabs / Courier rPr on the paragraph mark) untouched.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7Range = $p7.Range
$p7RunRange = $d.Range($p7Range.Start, $p7Range.End - 1)
$p7RunRange.InsertXML((New-WordOpenXml '<w:body><w:p><w:r><w:rPr/></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------------
# Change 4: stop.py -> stoplab (first occurrence, with "for the last time")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "before using “stop.py” to stop the lab for the last time.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "before using “stoplab” to stop the lab for the last time.",
    2)

# ---------------------------------------------------------------------------
# Change 5: "./stop.py pubkey" paragraph (with _GoBack bookmark) ->
# "stoplab pubkey", bookmark removed, paragraph mark run props cleared.
# ---------------------------------------------------------------------------
$p44 = $d.Paragraphs(44)
$p44Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:ind w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr>' + $runRPr + '</w:rPr><w:t>stoplab pubkey</w:t></w:r></w:p></w:body>'
$p44.Range.InsertXML((New-WordOpenXml $p44Body))

# ---------------------------------------------------------------------------
# Change 6: final "If you modified..." paragraph: clear paragraph mark run
# properties (Helvetica -> none) and stop.py -> stoplab in the text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "before typing “./stop.py”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "before typing “stoplab”.",
    2)

$p46 = $d.Paragraphs(46)
$p46Text = $p46.Range.Text
$p46Text = $p46Text.Substring(0, $p46Text.Length - 1)
$p46Body = '<w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs="Helvetica"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">' + $p46Text + '</w:t></w:r></w:p></w:body>'
$p46.Range.InsertXML((New-WordOpenXml $p46Body))
